# Re-rank / reshuffle the female real-effort ranking table so that the
# underlying worker-ranking data used to drive the instructions/infobox
# lines up correctly (see commit: "add infobox to make instructions for
# ranking easier/less").
#
# This updates the ranking scores (column F), the associated race value
# (column G), and swaps a handful of rows' prolific id / name / rank-input
# (columns B, C, D) so that each participant's identity travels with the
# correct record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 7.39000491208574

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("D3").Value = "Colleen"
$ws.Range("F3").Value = 6.091303748649244
$ws.Range("G3").Value = "White"

$ws.Range("B4").Value = 19
$ws.Range("C4").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("D4").Value = "Jewel"
$ws.Range("F4").Value = 6.07244745832581
$ws.Range("G4").Value = "Black or African American"

$ws.Range("B5").Value = 34
$ws.Range("C5").Value = "5e96194b0a9fe909389e9f7b"
$ws.Range("D5").Value = "Tina"
$ws.Range("F5").Value = 5.476255900907384
$ws.Range("G5").Value = "White"

$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "60bd88b8fc436774352f53b9"
$ws.Range("D6").Value = "Annes"
$ws.Range("F6").Value = 5.199586145078674
$ws.Range("G6").Value = "Asian"

$ws.Range("F7").Value = 4.182434273974512
$ws.Range("F8").Value = 1.405244066240008
$ws.Range("F9").Value = 1.285989741820513
$ws.Range("F10").Value = 0.4301310047900727
$ws.Range("F11").Value = 0.3677338533072753

$ws.Range("B12").Value = 33
$ws.Range("C12").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("D12").Value = "Shaniek"
$ws.Range("F12").Value = 0.2397686769137523

$ws.Range("B13").Value = 30
$ws.Range("C13").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("D13").Value = "Shadaisia"
$ws.Range("F13").Value = 0.2362855029629706
